# Apply the changes described by the diff:
#  - clear the stray numeric/date style (s="1") that LibreOffice's writer
#    left on the Boolean cells in each of the three sheets (the cells stay
#    t="b", they just go back to the "Normal" style)
#  - update the saved cell selection on each sheet
#  - leave "WithTable_Duplicate" as the active sheet/tab (was "Tableless")

$wb = $excel.ActiveWorkbook

# --- WithTable ---------------------------------------------------------
$wsWithTable = $wb.Worksheets.Item("WithTable")
$wsWithTable.Range("D2:D5").Style = "Normal"
$wsWithTable.Range("E23").Select()

# --- Tableless -----------------------------------------------------------
$wsTableless = $wb.Worksheets.Item("Tableless")
$wsTableless.Range("D2:D5").Style = "Normal"
$wsTableless.Range("I19").Select()

# --- WithTable_Duplicate --------------------------------------------------
$wsDup = $wb.Worksheets.Item("WithTable_Duplicate")
$wsDup.Range("E5:E8").Style = "Normal"
$wsDup.Range("M11").Select()
